$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 5
$ws.Range("E5").Value = 65

# Row 10
$ws.Range("E10").Value = 193

# Row 12
$ws.Range("E12").Value = 206

# Row 16
$ws.Range("E16").Value = 83
$ws.Range("F16").Value = 36
$ws.Range("H16").Value = 36

# Row 23
$ws.Range("E23").Value = 94

# Row 24
$ws.Range("F24").Value = 44
$ws.Range("H24").Value = 44

# Row 25
$ws.Range("E25").Value = 88
$ws.Range("F25").Value = 33
$ws.Range("H25").Value = 33

# Row 26
$ws.Range("E26").Value = 50

# Row 28
$ws.Range("E28").Value = 86
$ws.Range("F28").Value = 22
$ws.Range("H28").Value = 22

# Row 29
$ws.Range("E29").Value = 81

# Row 30
$ws.Range("E30").Value = 93
$ws.Range("F30").Value = 43
$ws.Range("H30").Value = 43

# Row 32
$ws.Range("E32").Value = 90

# Row 33
$ws.Range("E33").Value = 118

# Row 37
$ws.Range("F37").Value = 27
$ws.Range("H37").Value = 27

# Row 39
$ws.Range("E39").Value = 100

# Row 41
$ws.Range("E41").Value = 165

# Row 42
$ws.Range("E42").Value = 146

# Row 43
$ws.Range("E43").Value = 44

# Row 45
$ws.Range("E45").Value = 51
$ws.Range("F45").Value = 29
$ws.Range("H45").Value = 29

# Row 46
$ws.Range("E46").Value = 112
